# Append a new data row (row 87) to Sheet1, extending the used range from
# A1:K86 to A1:K87, matching the pattern of the existing rows (frn_adminid
# plus question1..question10 answers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 87
$values = @(197, 7, 4, 4, 5, 7, 7, 6, 5, 4, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $values[$i]
}
